# Bulk update student details: merge Sheet2 (roster) rows into Sheet1, add a
# new "stu_semester" column, and correct the one pre-existing student
# record's parent names.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- 1. Correct the pre-existing student's parent names first, while the
#        record is still sitting in its original row 2.
$ws1.Range("F2").Value2 = "Hamid Pasha"
$ws1.Range("H2").Value2 = "Sabnam Hamid"

# --- 2. Move that corrected record down to row 13 (last row), making room
#        to copy in the Sheet2 roster starting at row 2.
$ws1.Range("A13").Value2 = $ws1.Range("A2").Value2
$ws1.Range("B13").Value2 = $ws1.Range("B2").Value2
$ws1.Cells.Item(13,3).Formula = "=A13&`" `"&B13"
$ws1.Range("D13").Value2 = "I"
$ws1.Range("E13").Value2 = $ws1.Range("E2").Value2
$ws1.Range("F13").Value2 = $ws1.Range("F2").Value2
$ws1.Range("G13").Value2 = $ws1.Range("G2").Value2
$ws1.Range("H13").Value2 = $ws1.Range("H2").Value2
$ws1.Range("I13").Value2 = $ws1.Range("I2").Value2
$ws1.Range("J13").Value2 = $ws1.Range("J2").Value2
$ws1.Range("K13").Value2 = $ws1.Range("K2").Value2
$ws1.Range("L13").Value2 = $ws1.Range("L2").Value2

# --- 3. Header row: column D is now "stu_semester" (was "stu_class").
$ws1.Range("D1").Value2 = "stu_semester"

# --- 4. Copy the 11 student records from Sheet2 into Sheet1 rows 2-12,
#        inserting the constant semester value "I" into column D.
for ($r = 1; $r -le 11; $r++) {
    $dest = $r + 1

    $ws1.Cells.Item($dest,1).Value2 = $ws2.Cells.Item($r,1).Value2
    $ws1.Cells.Item($dest,2).Value2 = $ws2.Cells.Item($r,2).Value2
    $ws1.Cells.Item($dest,3).Formula = "=A$dest&`" `"&B$dest"
    $ws1.Cells.Item($dest,4).Value2 = "I"
    $ws1.Cells.Item($dest,5).Value2 = $ws2.Cells.Item($r,5).Value2
    $ws1.Cells.Item($dest,6).Value2 = $ws2.Cells.Item($r,6).Value2
    $ws1.Cells.Item($dest,7).Value2 = $ws2.Cells.Item($r,7).Value2
    $ws1.Cells.Item($dest,8).Value2 = $ws2.Cells.Item($r,8).Value2
    $ws1.Cells.Item($dest,9).Value2 = $ws2.Cells.Item($r,9).Value2
    $ws1.Cells.Item($dest,10).Value2 = $ws2.Cells.Item($r,10).Value2
    $ws1.Cells.Item($dest,11).Value2 = $ws2.Cells.Item($r,11).Value2
    $ws1.Cells.Item($dest,12).Value2 = $ws2.Cells.Item($r,12).Value2
}

# --- 5. Selection / view tidy-up to match the post-edit state. Select
#        Sheet2 first, then Sheet1 last, so Sheet1 ends up the active tab
#        again while both sheets keep their own remembered selection.
$ws2.Range("B32").Select()
$ws1.Range("D2").Select()
